# "New graphs for paper" - add two new mini-tables (error-vs-index and
# missing-data-rate series) to Sheet2 for the new paper figures, and widen /
# add the columns that hold them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Header row (row 2) ------------------------------------------------
# Order matters: new shared strings get interned in first-use order, and
# "Angle error" (Q2) must land before "Missing Data Rate %" (P2) so the
# shared-string table ends up Missing Data Rate (%), Standard Deviation...,
# Angle error, Missing Data Rate % (indices 107-110).
$ws.Range("J2").Value = "Missing Data Rate (%)"
$ws.Range("K2").Value = "Standard Deviation of Approximation Error Measurments"
$ws.Range("M2").Value = "Missing Data Rate (%)"
$ws.Range("Q2").Value = "Angle error"
$ws.Range("P2").Value = "Missing Data Rate %"

# --- Data rows (3-12): index columns J/M, error series K/N/Q ----------
$rows = @(
    @{ Row = 3;  J = 1;  K = 0.113078736087662;    N = 0.042039797521405999;  Q = 0.113078736087662 },
    @{ Row = 4;  J = 2;  K = 0.13565881390919601;  N = 0.057580257926705;     Q = 0.13565881390919601 },
    @{ Row = 5;  J = 3;  K = 0.15592321466704101;  N = 0.066301403203123693;  Q = 0.15592321466704101 },
    @{ Row = 6;  J = 4;  K = 0.17173736161074399;  N = 0.0537841045228254;    Q = 0.17173736161074399 },
    @{ Row = 7;  J = 5;  K = 0.156188257688813;    N = 0.069592245651760795;  Q = 0.156188257688813 },
    @{ Row = 8;  J = 6;  K = 0.163128550795407;    N = 0.0834708281066879;    Q = 0.163128550795407 },
    @{ Row = 9;  J = 7;  K = 0.207387696126413;    N = 0.083833530230314701;  Q = 0.207387696126413 },
    @{ Row = 10; J = 8;  K = 0.45668498876458802;  N = 0.15468930185714899;   Q = 0.45668498876458802 },
    @{ Row = 11; J = 9;  K = 0.43974896833390598;  N = 0.16920721323379501;   Q = 0.43974896833390598 },
    @{ Row = 12; J = 10; K = 0.28274015118616802;  N = 0.62769374040992398;   Q = 0.28274015118616802 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("M$row").Value = $r.J
    $ws.Range("N$row").Value = $r.N
    $ws.Range("Q$row").Value = $r.Q
}

# --- Column widths for the new / widened columns -----------------------
$ws.Columns.Item(10).ColumnWidth = 18.666666666666668   # J -> 19.5
$ws.Columns.Item(11).ColumnWidth = 51.833333333333336   # K -> 52.6640625 (widened from 43.6640625)
$ws.Columns.Item(13).ColumnWidth = 24.0                 # M -> 24.83203125
$ws.Columns.Item(16).ColumnWidth = 28.5                 # P -> 29.33203125

# --- View: new selection on the newly added data ------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 13
$ws.Range("P3").Select()
